$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "64.719.51"
$ws.Range("E2").Value = "  -2.24%  "

# Row 3
$ws.Range("D3").Value = "3.228.95"
$ws.Range("E3").Value = "  -1.44%  "

# Row 4
$ws.Range("E4").Value = "  +0.01%  "

# Row 5
$ws.Range("D5").Value = "'576.71"
$ws.Range("E5").Value = "  -1.55%  "

# Row 6
$ws.Range("D6").Value = "'172.06"
$ws.Range("E6").Value = "  -3.89%  "

# Row 7
$ws.Range("D7").Value = "'0.625"
$ws.Range("E7").Value = "  -0.22%  "

# Row 8
$ws.Range("E8").Value = "  +0.01%  "

# Row 9
$ws.Range("D9").Value = "3.226.55"
$ws.Range("E9").Value = "  -1.43%  "

# Row 10
$ws.Range("E10").Value = "  -2.83%  "

# Row 11
$ws.Range("D11").Value = "'6.76"
$ws.Range("E11").Value = "  +0.35%  "

# Row 12
$ws.Range("E12").Value = "  -3.10%  "

# Row 13
$ws.Range("D13").Value = "3.788.37"
$ws.Range("E13").Value = "  -1.52%  "

# Row 15
$ws.Range("D15").Value = "64.797.88"
$ws.Range("E15").Value = "  -2.12%  "

# Row 16
$ws.Range("D16").Value = "'25.73"
$ws.Range("E16").Value = "  -2.27%  "

# Row 17
$ws.Range("B17").Value = "ShibaInu"
$ws.Range("C17").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D17").Value = "'0.0000159"
$ws.Range("E17").Value = "  -3.03%  "

# Row 18
$ws.Range("B18").Value = "WrappedEther"
$ws.Range("C18").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D18").Value = "3.217.89"
$ws.Range("E18").Value = "  -1.98%  "

# Row 19
$ws.Range("D19").Value = "'416.29"
$ws.Range("E19").Value = "  -3.82%  "

# Row 20
$ws.Range("D20").Value = "'5.37"
$ws.Range("E20").Value = "  -2.40%  "

# Row 21
$ws.Range("D21").Value = "'12.82"
$ws.Range("E21").Value = "  -2.56%  "

# Row 22
$ws.Range("D22").Value = "'7.20"
$ws.Range("E22").Value = "  -2.38%  "

# Row 23
$ws.Range("D23").Value = "'0.999"
$ws.Range("E23").Value = "  -0.04%  "

# Row 24
$ws.Range("D24").Value = "'70.34"
$ws.Range("E24").Value = "  -1.80%  "

# Row 25
$ws.Range("E25").Value = "  -0.69%  "

# Row 26
$ws.Range("E26").Value = "  +4.30%  "

# Row 27
$ws.Range("D27").Value = "'0.495"
$ws.Range("E27").Value = "  -1.93%  "

# Row 28
$ws.Range("E28").Value = "  -1.86%  "

# Row 29
$ws.Range("D29").Value = "'8.96"
$ws.Range("E29").Value = "  +1.55%  "

# Row 30
$ws.Range("E30").Value = "  +0.14%  "

# Row 31
$ws.Range("D31").Value = "'1.86"
$ws.Range("E31").Value = "  -4.73%  "

# Row 32
$ws.Range("D32").Value = "'21.77"
$ws.Range("E32").Value = "  -2.08%  "

# Row 33
$ws.Range("D33").Value = "'0.999"
$ws.Range("E33").Value = "  +0.05%  "

# Row 34
$ws.Range("E34").Value = "  -3.18%  "

# Row 35
$ws.Range("D35").Value = "'6.41"

# Row 36
$ws.Range("D36").Value = "'1.15"
$ws.Range("E36").Value = "  -2.78%  "

# Row 37
$ws.Range("D37").Value = "'157.91"
$ws.Range("E37").Value = "  +0.04%  "

# Row 38
$ws.Range("E38").Value = "  -1.83%  "

# Row 39
$ws.Range("D39").Value = "2.809.96"
$ws.Range("E39").Value = "  +1.27%  "

# Row 40
$ws.Range("E40").Value = "  -2.93%  "

# Row 41
$ws.Range("D41").Value = "'25.41"
$ws.Range("E41").Value = "  -4.04%  "

# Row 42
$ws.Range("E42").Value = "  -2.69%  "

# Row 43
$ws.Range("D43").Value = "'39.41"
$ws.Range("E43").Value = "  -1.93%  "

# Row 44
$ws.Range("D44").Value = "'0.722"
$ws.Range("E44").Value = "  -6.48%  "

# Row 45
$ws.Range("D45").Value = "'5.76"
$ws.Range("E45").Value = "  -4.25%  "

# Row 46
$ws.Range("D46").Value = "'0.0628"
$ws.Range("E46").Value = "  -4.43%  "

# Row 47
$ws.Range("D47").Value = "'2.18"
$ws.Range("E47").Value = "  -4.49%  "

# Row 48
$ws.Range("D48").Value = "'302.00"
$ws.Range("E48").Value = "  -5.88%  "

# Row 49
$ws.Range("D49").Value = "'21.97"
$ws.Range("E49").Value = "  -5.21%  "

# Row 50
$ws.Range("E50").Value = "  -1.33%  "

# Row 51
$ws.Range("E51").Value = "  -1.69%  "
